$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 25.16000000000049
$ws.Range("H2").Value = 0.4223512403010066
$ws.Range("I2").Value = 0.4223512403010066
$ws.Range("L2").Value = 4.42269532729915
$ws.Range("M2").Value = "[-4.439608803572723, 13.284999458171024]"
$ws.Range("N2").Value = 0.3202087748901186
$ws.Range("O2").Value = 0.3202087748901186
$ws.Range("P2").Value = -1.786210838077079
$ws.Range("Q2").Value = "[-4.9183692794939295, 1.3459476033397708]"
$ws.Range("R2").Value = 0.2567878470501963
$ws.Range("S2").Value = 0.2567878470501963
$ws.Range("T2").Value = 14.85124683973271
$ws.Range("U2").Value = "[10.018168205488465, 19.68432547397695]"
$ws.Range("V2").Value = 0.0000001635629234275626
$ws.Range("W2").Value = 0.0000001635629234275626
$ws.Range("X2").Value = 7.152592592592733
$ws.Range("Y2").Value = -5.389629629629734
$ws.Range("Z2").Value = 19.6948148148152
$ws.Range("F3").Value = 25.16000000000049
$ws.Range("H3").Value = 0.5964868600298689
$ws.Range("I3").Value = 0.5964868600298689
$ws.Range("L3").Value = 3.851235008294919
$ws.Range("M3").Value = "[-5.956207114232749, 13.658677130822587]"
$ws.Range("N3").Value = 0.4331464220726478
$ws.Range("O3").Value = 0.4331464220726478
$ws.Range("P3").Value = 2.886868924920351
$ws.Range("Q3").Value = "[-0.23900004171453837, 6.01273789155524]"
$ws.Range("R3").Value = 0.06941377427731221
$ws.Range("S3").Value = 0.06941377427731221
$ws.Range("T3").Value = 14.98009383746181
$ws.Range("U3").Value = "[9.690129660441528, 20.270058014482085]"
$ws.Range("V3").Value = 0.0000008584345276840821
$ws.Range("W3").Value = 0.0000008584345276840821
$ws.Range("X3").Value = 13.60000000000027
$ws.Range("Y3").Value = 1.082962962962984
$ws.Range("Z3").Value = 26.11703703703755
$ws.Range("B4").Value = 0
$ws.Range("F4").Value = 25.16000000000049
$ws.Range("H4").Value = 0.1217080612877719
$ws.Range("I4").Value = 0.1217080612877719
$ws.Range("L4").Value = 7.273897311667289
$ws.Range("M4").Value = "[-2.1025600027723534, 16.65035462610693]"
$ws.Range("N4").Value = 0.1251846183985066
$ws.Range("O4").Value = 0.1251846183985066
$ws.Range("P4").Value = -2.44031621540108
$ws.Range("Q4").Value = "[-5.434106211614777, 0.5534737808126167]"
$ws.Range("R4").Value = 0.1076136729251054
$ws.Range("S4").Value = 0.1076136729251054
$ws.Range("T4").Value = 16.71124804867774
$ws.Range("U4").Value = "[11.84794182499669, 21.57455427235879]"
$ws.Range("V4").Value = 0.00000001335347632469563
$ws.Range("W4").Value = 0.00000001335347632469563
$ws.Range("X4").Value = 9.771851851852043
$ws.Range("Y4").Value = -2.216296296296342
$ws.Range("Z4").Value = 21.76000000000043
$ws.Range("F5").Value = 25.16000000000049
$ws.Range("H5").Value = 0.2661205318205228
$ws.Range("I5").Value = 0.2661205318205228
$ws.Range("L5").Value = 6.951451267702248
$ws.Range("M5").Value = "[-4.733784067249405, 18.6366866026539]"
$ws.Range("N5").Value = 0.2371219786907606
$ws.Range("O5").Value = 0.2371219786907606
$ws.Range("P5").Value = -2.377421467581465
$ws.Range("Q5").Value = "[-5.515869383780277, 0.7610264486173479]"
$ws.Range("R5").Value = 0.1340795459855699
$ws.Range("S5").Value = 0.1340795459855699
$ws.Range("T5").Value = 15.52829964003021
$ws.Range("U5").Value = "[9.59479619290498, 21.46180308715543]"
$ws.Range("V5").Value = 0.000003712527029353652
$ws.Range("W5").Value = 0.000003712527029353652
$ws.Range("X5").Value = 9.520000000000184
$ws.Range("Y5").Value = -3.047407407407469
$ws.Range("Z5").Value = 22.08740740740784
$ws.Range("F6").Value = 25.16000000000049
$ws.Range("H6").Value = 0.3068419013648915
$ws.Range("I6").Value = 0.3068419013648915
$ws.Range("L6").Value = 5.560711984506323
$ws.Range("M6").Value = "[-3.236258504564997, 14.357682473577643]"
$ws.Range("N6").Value = 0.2095033462294329
$ws.Range("O6").Value = 0.2095033462294329
$ws.Range("P6").Value = 3.050395269251351
$ws.Range("Q6").Value = "[-0.08805264694746207, 6.188843185450164]"
$ws.Range("R6").Value = 0.05649643024470974
$ws.Range("S6").Value = 0.05649643024470974
$ws.Range("T6").Value = 14.28521430121324
$ws.Range("U6").Value = "[9.302423343373144, 19.268005259053346]"
$ws.Range("V6").Value = 0.0000006746996643602898
$ws.Range("W6").Value = 0.0000006746996643602898
$ws.Range("X6").Value = 12.94518518518544
$ws.Range("Y6").Value = 0.3777777777777818
$ws.Range("Z6").Value = 25.51259259259309
$ws.Range("F7").Value = 25.03000000000047
$ws.Range("H7").Value = 0.5046369962220019
$ws.Range("I7").Value = 0.5046369962220019
$ws.Range("L7").Value = 4.142840444777678
$ws.Range("M7").Value = "[-4.887098585447927, 13.172779475003285]"
$ws.Range("N7").Value = 0.360389312464283
$ws.Range("O7").Value = 0.360389312464283
$ws.Range("P7").Value = 2.937184723176043
$ws.Range("Q7").Value = "[-0.19497371824080822, 6.069343164592894]"
$ws.Range("R7").Value = 0.06538390748816458
$ws.Range("S7").Value = 0.06538390748816458
$ws.Range("T7").Value = 14.44024535837705
$ws.Range("U7").Value = "[9.510696819315454, 19.36979389743864]"
$ws.Range("V7").Value = 0.0000004394402821628063
$ws.Range("W7").Value = 0.0000004394402821628063
$ws.Range("X7").Value = 13.32928928928954
$ws.Range("Y7").Value = 0.8518718718718876
$ws.Range("Z7").Value = 25.8067067067072
$ws.Range("F8").Value = 25.03000000000047
$ws.Range("H8").Value = 0.2347771966479127
$ws.Range("I8").Value = 0.2347771966479127
$ws.Range("L8").Value = 5.690092781486698
$ws.Range("M8").Value = "[-3.079163905217678, 14.459349468191073]"
$ws.Range("N8").Value = 0.1978901576451961
$ws.Range("O8").Value = 0.1978901576451961
$ws.Range("P8").Value = 2.55981623625835
$ws.Range("Q8").Value = "[-0.5786316799404627, 5.698264152457163]"
$ws.Range("R8").Value = 0.1074019953907757
$ws.Range("S8").Value = 0.1074019953907757
$ws.Range("T8").Value = 14.05554119568906
$ws.Range("U8").Value = "[9.37400667635087, 18.737075715027256]"
$ws.Range("V8").Value = 0.0000002658750848372193
$ws.Range("W8").Value = 0.0000002658750848372193
$ws.Range("X8").Value = 14.83259259259287
$ws.Range("Y8").Value = 2.33012012012016
$ws.Range("Z8").Value = 27.33506506506559
$ws.Range("F9").Value = 25.03000000000047
$ws.Range("H9").Value = 0.1771226070754988
$ws.Range("I9").Value = 0.1771226070754988
$ws.Range("L9").Value = 6.246547452384466
$ws.Range("M9").Value = "[-2.3474212119556297, 14.840516116724562]"
$ws.Range("N9").Value = 0.1501570904080856
$ws.Range("O9").Value = 0.1501570904080856
$ws.Range("P9").Value = 2.673026782333658
$ws.Range("Q9").Value = "[-0.30818426431611634, 5.654237828983431]"
$ws.Range("R9").Value = 0.07762627987293391
$ws.Range("S9").Value = 0.07762627987293391
$ws.Range("T9").Value = 13.91284364808741
$ws.Range("U9").Value = "[9.246168273038439, 18.57951902313639]"
$ws.Range("V9").Value = 0.0000003072800065506698
$ws.Range("W9").Value = 0.0000003072800065506698
$ws.Range("X9").Value = 14.38160160160188
$ws.Range("Y9").Value = 2.505505505505553
$ws.Range("Z9").Value = 26.2576976976982
$ws.Range("F10").Value = 25.03000000000047
$ws.Range("H10").Value = 0.05054397120082244
$ws.Range("I10").Value = 0.05054397120082244
$ws.Range("L10").Value = 7.760919992331388
$ws.Range("M10").Value = "[-0.7294699985186668, 16.251309983181443]"
$ws.Range("N10").Value = 0.07221145771690218
$ws.Range("O10").Value = 0.07221145771690218
$ws.Range("P10").Value = 2.245342497160272
$ws.Range("Q10").Value = "[0.018868424345884982, 4.471816569974659]"
$ws.Range("R10").Value = 0.04816948331393522
$ws.Range("S10").Value = 0.04816948331393522
$ws.Range("T10").Value = 11.09751488413155
$ws.Range("U10").Value = "[6.66328008284932, 15.531749685413786]"
$ws.Range("V10").Value = 0.000008032878042074643
$ws.Range("W10").Value = 0.000008032878042074643
$ws.Range("X10").Value = 16.08534534534565
$ws.Range("Y10").Value = 7.215855855855997
$ws.Range("Z10").Value = 24.9548348348353
$ws.Range("F11").Value = 25.03000000000047
$ws.Range("H11").Value = 0.4895275894535767
$ws.Range("I11").Value = 0.4895275894535767
$ws.Range("L11").Value = 4.538408009504867
$ws.Range("M11").Value = "[-5.710140866696143, 14.786956885705877]"
$ws.Range("N11").Value = 0.3771837026959259
$ws.Range("O11").Value = 0.3771837026959259
$ws.Range("P11").Value = 2.044079304137503
$ws.Range("Q11").Value = "[-1.0817896624973864, 5.169948270772392]"
$ws.Range("R11").Value = 0.1944829485704491
$ws.Range("S11").Value = 0.1944829485704491
$ws.Range("T11").Value = 13.17527163013316
$ws.Range("U11").Value = "[7.710125533876648, 18.64041772638967]"
$ws.Range("V11").Value = 0.00001485673553602851
$ws.Range("W11").Value = 0.00001485673553602851
$ws.Range("X11").Value = 16.88710710710743
$ws.Range("Y11").Value = 4.434744744744828
$ws.Range("Z11").Value = 29.33946946947002
$ws.Range("F12").Value = 25.03000000000047
$ws.Range("H12").Value = 0.1194584660384148
$ws.Range("I12").Value = 0.1194584660384148
$ws.Range("L12").Value = 6.138747496983511
$ws.Range("M12").Value = "[-1.2967168582472883, 13.57421185221431]"
$ws.Range("N12").Value = 0.1032929623548693
$ws.Range("O12").Value = 0.1032929623548693
$ws.Range("P12").Value = 1.452868674633118
$ws.Range("Q12").Value = "[-0.1698158191129604, 3.075553168379196]"
$ws.Range("R12").Value = 0.07803499052242402
$ws.Range("S12").Value = 0.07803499052242402
$ws.Range("T12").Value = 12.53987418312434
$ws.Range("U12").Value = "[8.3224879620157, 16.757260404232976]"
$ws.Range("V12").Value = 0.00000032454814258287
$ws.Range("W12").Value = 0.00000032454814258287
$ws.Range("X12").Value = 19.24228228228264
$ws.Range("Y12").Value = 12.77807807807832
$ws.Range("Z12").Value = 25.70648648648697
